$wb = $excel.ActiveWorkbook

# --- Insert the new "registrationSheet" right after "loginSheet" ---
$loginSheet = $wb.Worksheets.Item("loginSheet")
$regSheet = $wb.Worksheets.Add($null, $loginSheet)
$regSheet.Name = "registrationSheet"

# --- Column headers (username/password reuse existing shared strings) ---
$regSheet.Range("A1").Value = "username"
$regSheet.Range("B1").Value = "password"

# --- Registration test data, entered row by row ---
$regSheet.Range("A3").Value = "NumpyTest1"

$regSheet.Range("B4").Value = "Testpassword2"

$regSheet.Range("A5").Value = "abcd%"
$regSheet.Range("D5").Value = "Please enter a valid username"

$regSheet.Range("A6").Value = "Numpyninja"
$regSheet.Range("B6").Value = "lessnum"
$regSheet.Range("D6").Value = "Password should contain at least 8 characters"

# --- Remaining header cells ---
$regSheet.Range("C1").Value = "confirmPassword"
$regSheet.Range("D1").Value = "expectedMessage"
$regSheet.Range("D2").Value = "Please fill out this field."

# --- Validation tag column ---
$regSheet.Range("E3").Value = "password_empty"
$regSheet.Range("E6").Value = "password_invalid"
$regSheet.Range("E1").Value = "validation"
$regSheet.Range("E2").Value = "username_empty"
$regSheet.Range("E4").Value = "confirmpassword_empty"
$regSheet.Range("E5").Value = "username_invalid"

# --- Fill in the remaining (repeated) values ---
$regSheet.Range("A4").Value = "NumpyTest1"
$regSheet.Range("B5").Value = "Testpassword2"
$regSheet.Range("C5").Value = "Testpassword2"
$regSheet.Range("C6").Value = "lessnum"
$regSheet.Range("D3").Value = "Please fill out this field."
$regSheet.Range("D4").Value = "Please fill out this field."

# --- Formatting: wrap + vertically centered body cells ---
$regSheet.Range("A2:D6").WrapText = $true
$regSheet.Range("A2:D6").VerticalAlignment = -4108

# --- Italicize the trailing "%" of the invalid-username sample ---
$regSheet.Range("A5").Characters(5, 1).Font.Italic = $true
$regSheet.Range("A5").Font.Italic = $true

# --- Column widths to match the new content ---
$regSheet.Columns("A").ColumnWidth = 15.88671875
$regSheet.Columns("B").ColumnWidth = 15.44140625
$regSheet.Columns("C").ColumnWidth = 16.5546875
$regSheet.Columns("D").ColumnWidth = 41.21875
$regSheet.Columns("E").ColumnWidth = 35.44140625

# --- Selection state: loginSheet keeps header row selected, registrationSheet becomes the active tab ---
$loginSheet.Rows("1:1").Select() | Out-Null

$regSheet.Range("E5").Select() | Out-Null
